$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column K (value + same header formatting as the other header cells)
$ws.Range("K1").Value = "intervention_type"
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)  # xlPasteFormats

# Values for new column K, rows 2-20
$values = @(
    "OTHER",      # row 2
    "BEHAVIORAL", # row 3
    "BEHAVIORAL", # row 4
    "PROCEDURE",  # row 5
    "BEHAVIORAL", # row 6
    "DRUG",       # row 7
    "BEHAVIORAL", # row 8
    "DEVICE",     # row 9
    "OTHER",      # row 10
    "PROCEDURE",  # row 11
    "OTHER",      # row 12
    "OTHER",      # row 13
    "DRUG",       # row 14
    "OTHER",      # row 15
    "OTHER",      # row 16
    "OTHER",      # row 17
    "BIOLOGICAL", # row 18
    "OTHER",      # row 19
    "OTHER"       # row 20
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 11).Value = $values[$i]
}
